# This edit re-orders the data rows (rows 2-9) of the sheet: the records are
# cyclically shuffled into a different row order while all other structure
# (header row, columns, formatting) stays the same.
#
# Mapping of new row -> old row it is populated from (1-based worksheet rows):
#   2 <- 7
#   3 <- 2
#   4 <- 3
#   5 <- 8
#   6 <- 4
#   7 <- 9
#   8 <- 5
#   9 <- 6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 9
$lastCol  = 51   # column AY

# A couple of columns (Startdatum / Slutdatum) contain plain text values that
# look like dates ("2018-10-01"). Force those columns to stay text so Excel
# does not silently reinterpret them as date serial numbers when the values
# are written back.
$ws.Range("Y$($firstRow):Y$($lastRow)").NumberFormat = "@"
$ws.Range("AA$($firstRow):AA$($lastRow)").NumberFormat = "@"

# Snapshot the current contents of the data rows before overwriting anything.
$sourceRange = $ws.Range("A$($firstRow):AY$($lastRow)")
$original = $sourceRange.Value2

$rowCount = $lastRow - $firstRow + 1

# new row index (1-based, relative to $firstRow) -> old row index (same basis)
$mapping = @{
    1 = 6   # row 2 <- row 7
    2 = 1   # row 3 <- row 2
    3 = 2   # row 4 <- row 3
    4 = 7   # row 5 <- row 8
    5 = 3   # row 6 <- row 4
    6 = 8   # row 7 <- row 9
    7 = 4   # row 8 <- row 5
    8 = 5   # row 9 <- row 6
}

$updated = New-Object 'object[,]' $rowCount, $lastCol
for ($newIdx = 1; $newIdx -le $rowCount; $newIdx++) {
    $oldIdx = $mapping[$newIdx]
    for ($c = 1; $c -le $lastCol; $c++) {
        $updated[$newIdx - 1, $c - 1] = $original[$oldIdx, $c]
    }
}

$ws.Range("A$($firstRow):AY$($lastRow)").Value2 = $updated
